$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 (2021Q2)
$ws.Range("C4").Value = 188
$ws.Range("E4").Value = 107

# Row 5 (2021Q3)
$ws.Range("C5").Value = 215
$ws.Range("D5").Value = 96
$ws.Range("E5").Value = 119
$ws.Range("F5").Value = 51.06382978723404

# Row 6 (2021Q4)
$ws.Range("C6").Value = 217
$ws.Range("D6").Value = 129
$ws.Range("F6").Value = 60

# Row 7 (2022Q1)
$ws.Range("C7").Value = 237
$ws.Range("D7").Value = 130
$ws.Range("E7").Value = 107
$ws.Range("F7").Value = 59.90783410138248

# Row 8 (2022Q2)
$ws.Range("F8").Value = 54.85232067510548

# Row 9 (2022Q3)
$ws.Range("C9").Value = 177
$ws.Range("E9").Value = 51

# Row 10 (2022Q4)
$ws.Range("C10").Value = 165
$ws.Range("D10").Value = 125
$ws.Range("F10").Value = 70.62146892655367

# Row 11 (2023Q1)
$ws.Range("C11").Value = 195
$ws.Range("D11").Value = 143
$ws.Range("F11").Value = 86.66666666666667

# Row 12 (2023Q2)
$ws.Range("C12").Value = 201
$ws.Range("D12").Value = 153
$ws.Range("F12").Value = 78.46153846153847

# Row 13 (2023Q3)
$ws.Range("C13").Value = 215
$ws.Range("D13").Value = 159
$ws.Range("F13").Value = 79.1044776119403

# Row 14 (2023Q4)
$ws.Range("C14").Value = 225
$ws.Range("D14").Value = 183
$ws.Range("E14").Value = 42
$ws.Range("F14").Value = 85.11627906976744

# Row 15 (2024Q1)
$ws.Range("C15").Value = 255
$ws.Range("D15").Value = 186
$ws.Range("E15").Value = 69
$ws.Range("F15").Value = 82.66666666666667

# Row 16 (2024Q2)
$ws.Range("C16").Value = 313
$ws.Range("D16").Value = 209
$ws.Range("E16").Value = 104
$ws.Range("F16").Value = 81.96078431372548

# Row 17 (2024Q3)
$ws.Range("C17").Value = 317
$ws.Range("D17").Value = 250
$ws.Range("E17").Value = 67
$ws.Range("F17").Value = 79.87220447284345

# Row 18 (2024Q4)
$ws.Range("C18").Value = 319
$ws.Range("D18").Value = 255
$ws.Range("E18").Value = 64
$ws.Range("F18").Value = 80.4416403785489

# Row 19 (2025Q1)
$ws.Range("C19").Value = 323
$ws.Range("D19").Value = 263
$ws.Range("E19").Value = 60
$ws.Range("F19").Value = 82.44514106583071

# Row 20 (2025Q2)
$ws.Range("C20").Value = 319
$ws.Range("D20").Value = 256
$ws.Range("E20").Value = 63
$ws.Range("F20").Value = 79.25696594427245
